$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New (encrypted-looking) payload strings that replace the previous
# shared-string pair used in column B.
$secretC = 'U2FsdGVkX1+SHshU4/0vmkrFjgM5GaPFhmwH220YLrbOM3uP+towDO/6RJ4AqiX/f9uqXP8kR3O2d8Zkig9+JjD9w90DNg3soxrSShkFADlEnhB4wTrx8Gp6XzTDDAswJFWGfk+Qb9/QiCJLCHfYfIXM8B1lk9Fc7pieqDEDLYc9z0q9Nwmn5g8c1TQuEsnvsykRXIOF/is9tdEdboCsYO5HoozfzjJxQ9pw9RFwQHxI8O8UkFVQ5LbEF1jvzhpTU2sD0t5v2qnRK5tag95Qs/IpmWpCDevexzdGYY+VrM8Yvmk8KWekLiW4MCmYgERQIbTI2nXNrGa9tIceAk+f0uZN9gEgdZH72R0tkhBBWx7p+DPk8GQMsVbK3hC1UJBJ8Em8vpkM6O7jED/QnayKYM17/o0FUNK/yozpvmqsb9AJdhcU0ZqlLEoLRPSG/iPZh2GD2upUXf4RuhrgpywGsfRuwgKWKdLEXSwg9RQFU8Yo1nctY//HHcBMU+GORH3gkHjC2n3i+f/us6XCX08GOvArx2/aYfyfndw5uxv38egarzYAWnTZ+5Iaetuu0Wz7qBbZqZN5fXlQ8PaVel8bWUGRKm8VuKKIi6XR5EtgoRyEi6NaamOVgjWKC4MARK87aB1UcJX2K+nS3bsurIyvLD12wTOwpo8d5EwuZ1tsBTrHE0PRmEMlJjVIn5KgNOn04NVpA4eTOhYvhlTLJzUBhH47iLnTvRbTmLto0dcvNKuVgXiTcG1we7AtX871h6WyeBA79WiYh6CY3hPusALKJUI1JGjWDCmNdVKfk5RtljmdtW+WCJxVAl6fUjyqIwMzISoiv+PkrjS6Oj9zX7bZwtDfzGRVYQaw4Ur22/tA1fpwmq4C+jtfyZcinqB20U37ASCCXfCSu1ANlnfnumLjLJ6k6vTvchC8l1hPHprWr/TF91XKIBJJCe4Dc0uZgAGjOCmfOlXH2fCGbYYGJgqGYQaBozpD13Y5AyWqPRhZtSN169qxqOYERhDYDpTCfpGBbY8URcB8CMl92PBosttAGR0db8bZzIwkHSbu7BXG2K9OmFgwjgMvO0nBNH8OdK6sSFp+bfk0gCdmOp8wkSf5NFhudsTl9jH60ZZEa0RZSzd1POjirTwFqy/PVviuILEQccz67eCDlFKnx6/pX3uI4oV3c2oKXmj7BuXS8c5RR4SZeAQA/o1ny6gKdn4xRKuvBF287b32JTAZ2hM5wGuuojqx88VLtSr7+eIurJS8ga3ri0XV0Rb6J48k5NT9F7qKwZ96XHseCT18Xf0mYCNKiAd6w21BOnJaaunNWm3L93g='
$secretD = 'U2FsdGVkX1//PByWnO0jHrn2mZx5+i/56n9OWuHlT2nh176eIoNaiMC3x+pvJRYIywkvduW5FvX31S9wPxI8E9DBVasntf44SaGcagTmpxpvQm1XLe+bxXtUaDZhX8jjflZPhzuCICW3gO7HmAoMxD/MNAb1KzaevR/eQw06HGEXYQ+pT07loWptcw37RpXVC0LtUcLcgoL4QSyRoYU8nDG+NRDKMw/26oicMj+o7xjqlJr3kZI+OuaQKaVeVy0yoLbxCnuqcuzkfx5Rpg2hyxUmoZvScoiXsCj4a0SCAiIAqlsZowCopGowcBdefzuH'

# Column B alternated between two placeholder values row-by-row; replace
# each occurrence with the corresponding new payload, keeping the same
# alternating pattern and leaving column C (time ranges) untouched.
# Saving in between the two batches preserves the original relative
# ordering of the two replacement strings inside the shared-string table.
$ws.Range("B2").Value = $secretC
$ws.Range("B4").Value = $secretC
$ws.Range("B6").Value = $secretC
$wb.Save()

$ws.Range("B3").Value = $secretD
$ws.Range("B5").Value = $secretD
$ws.Range("B7").Value = $secretD

# Restore the active selection to match the saved view state.
$ws.Range("B10").Select() | Out-Null
